$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) / Volume(1h) (E) values scraped on 2023-07-23.
# Key = row number; Price = $null when unchanged this run.
$updates = @(
    @{ Row = 2; Price = "29.900.10"; Volume = "  -0.24%  " }
    @{ Row = 3; Price = "1.873.79"; Volume = "  -1.06%  " }
    @{ Row = 4; Price = $null; Volume = "  +0.05%  " }
    @{ Row = 5; Price = "0.7378"; Volume = "  -4.83%  " }
    @{ Row = 6; Price = "242.33"; Volume = "  -0.66%  " }
    @{ Row = 7; Price = "1.000"; Volume = "  +0.00%  " }
    @{ Row = 8; Price = "0.3155"; Volume = "  +0.74%  " }
    @{ Row = 9; Price = "0.07174"; Volume = "  -1.19%  " }
    @{ Row = 10; Price = "24.67"; Volume = "  -4.45%  " }
    @{ Row = 11; Price = "0.08412"; Volume = "  -2.86%  " }
    @{ Row = 12; Price = "0.7501"; Volume = "  -3.01%  " }
    @{ Row = 13; Price = "5.415"; Volume = "  -0.10%  " }
    @{ Row = 14; Price = "1.877.18"; Volume = "  -8.78%  " }
    @{ Row = 15; Price = "92.52"; Volume = "  -2.11%  " }
    @{ Row = 16; Price = "29.909.33"; Volume = "  -0.78%  " }
    @{ Row = 17; Price = "6.098"; Volume = "  -1.81%  " }
    @{ Row = 18; Price = "13.58"; Volume = "  -2.69%  " }
    @{ Row = 19; Price = "242.70"; Volume = "  -1.27%  " }
    @{ Row = 20; Price = "0.000007811"; Volume = "  -1.06%  " }
    @{ Row = 22; Price = "2.117.30"; Volume = "  -8.33%  " }
    @{ Row = 23; Price = "7.999"; Volume = "  -2.14%  " }
    @{ Row = 24; Price = "1.000"; Volume = "  -0.04%  " }
    @{ Row = 25; Price = $null; Volume = "  -3.06%  " }
    @{ Row = 26; Price = "9.282"; Volume = "  -2.78%  " }
    @{ Row = 27; Price = "164.91"; Volume = "  +1.25%  " }
    @{ Row = 28; Price = "18.60"; Volume = "  -1.42%  " }
    @{ Row = 29; Price = "2.034"; Volume = "  -0.73%  " }
    @{ Row = 30; Price = "1.490"; Volume = "  +4.04%  " }
    @{ Row = 31; Price = "4.594"; Volume = "  +1.34%  " }
    @{ Row = 32; Price = "1.533"; Volume = "  -0.88%  " }
    @{ Row = 33; Price = "4.248"; Volume = "  +2.77%  " }
    @{ Row = 34; Price = "0.05318"; Volume = "  -2.44%  " }
    @{ Row = 35; Price = "1.234"; Volume = "  -1.35%  " }
    @{ Row = 36; Price = "0.7543"; Volume = "  -0.09%  " }
    @{ Row = 37; Price = "0.9982"; Volume = "  -0.24%  " }
    @{ Row = 38; Price = "2.691"; Volume = "  +0.19%  " }
    @{ Row = 39; Price = "0.01947"; Volume = "  -1.67%  " }
    @{ Row = 40; Price = "2.754"; Volume = "  -1.14%  " }
    @{ Row = 41; Price = "0.4509"; Volume = "  -0.36%  " }
    @{ Row = 42; Price = "1.110.53"; Volume = "  +1.67%  " }
    @{ Row = 43; Price = "6.053"; Volume = "  -0.15%  " }
    @{ Row = 44; Price = "72.16"; Volume = "  -1.92%  " }
    @{ Row = 45; Price = "0.8565"; Volume = "  +0.26%  " }
    @{ Row = 46; Price = "1.001"; Volume = "  +0.14%  " }
    @{ Row = 47; Price = "103.15"; Volume = "  -0.18%  " }
    @{ Row = 48; Price = "7.647"; Volume = "  +0.18%  " }
    @{ Row = 49; Price = "3.098"; Volume = "  +3.29%  " }
    @{ Row = 50; Price = "1.838"; Volume = "  -2.61%  " }
    @{ Row = 51; Price = "2.015.24"; Volume = "  -8.15%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $priceCell = $ws.Range("D" + $u.Row)
        # Force text so values like "1.000" / "29.900.10" are not reinterpreted as numbers.
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.Price
        $priceCell.Style = "Normal"
    }
    if ($null -ne $u.Volume) {
        $ws.Range("E" + $u.Row).Value = $u.Volume
    }
}
